# Insert a new row at row 257, shifting the existing rows 257:369 down to 258:370,
# then populate the newly inserted row with the new weekly price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("257:257").Insert()

$ws.Range("A257").Value = 10
$ws.Range("B257").Value = "Vega Modelo de Temuco"
$ws.Range("C257").Value = "La Araucanía"
$ws.Range("D257").Value = 44726
$ws.Range("E257").Value = 9
$ws.Range("F257").Value = "Fruta"
$ws.Range("G257").Value = 100108
$ws.Range("H257").Value = "Tropicales y subtropicales"
$ws.Range("I257").Value = 100108002
$ws.Range("J257").Value = "Mango"
$ws.Range("K257").Value = "Sin especificar"
$ws.Range("L257").Value = "Primera"
$ws.Range("M257").Value = 1200
$ws.Range("N257").Value = 9000
$ws.Range("O257").Value = 10000
$ws.Range("P257").Value = 9583
$ws.Range("Q257").Value = "$/bandeja 4 kilos"
$ws.Range("R257").Value = "Brasil"
$ws.Range("S257").Value = 2396
$ws.Range("T257").Value = 4
